$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 1.8
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 1.06
$ws.Range("K5").Value = 10
$ws.Range("L5").Value = 1.25
$ws.Range("M5").Value = 3.75
$ws.Range("N5").Value = 1.9
$ws.Range("O5").Value = 1.9
$ws.Range("U5").Value = 8.5
$ws.Range("AD5").Value = 13
$ws.Range("AE5").Value = 23

# Row 8 updates
$ws.Range("N8").Value = 1.92
$ws.Range("O8").Value = 1.82
